$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell text values are updated in place. For D-column price cells whose new
# text would otherwise be auto-parsed by Excel as a number (losing the exact
# original text formatting, e.g. trailing zeros or precision), we prefix the
# assignment with a leading apostrophe to force literal text, then restore the
# cells original style afterwards so no formatting/style diff is introduced.

$ws.Range("D2").Value = "26.510.62"
$ws.Range("E2").Value = "  +2.03%  "
$ws.Range("D3").Value = "1.681.95"
$ws.Range("E3").Value = "  +2.56%  "
$origStyle = $ws.Range("D4").Style
$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = $origStyle
$ws.Range("E4").Value = "  -0.09%  "
$origStyle = $ws.Range("D5").Style
$ws.Range("D5").Value = "'217.49"
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = "  +3.67%  "
$origStyle = $ws.Range("D6").Style
$ws.Range("D6").Value = "'0.5323"
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = "  +3.00%  "
$ws.Range("E7").Value = "  -0.08%  "
$origStyle = $ws.Range("D8").Style
$ws.Range("D8").Value = "'0.2674"
$ws.Range("D8").Style = $origStyle
$ws.Range("E8").Value = "  +4.30%  "
$origStyle = $ws.Range("D9").Style
$ws.Range("D9").Value = "'0.06426"
$ws.Range("D9").Style = $origStyle
$ws.Range("E9").Value = "  +3.06%  "
$origStyle = $ws.Range("D10").Style
$ws.Range("D10").Value = "'21.46"
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = "  +5.34%  "
$origStyle = $ws.Range("D11").Style
$ws.Range("D11").Value = "'0.07789"
$ws.Range("D11").Style = $origStyle
$ws.Range("E11").Value = "  +3.50%  "
$ws.Range("D12").Value = "1.691.29"
$ws.Range("E12").Value = "  +3.01%  "
$origStyle = $ws.Range("D13").Style
$ws.Range("D13").Value = "'4.509"
$ws.Range("D13").Style = $origStyle
$ws.Range("E13").Value = "  +3.46%  "
$origStyle = $ws.Range("D14").Style
$ws.Range("D14").Value = "'0.5623"
$ws.Range("D14").Style = $origStyle
$ws.Range("E14").Value = "  +4.08%  "
$ws.Range("D15").Value = "0.0₅8430"
$ws.Range("E15").Value = "  +5.97%  "
$origStyle = $ws.Range("D16").Style
$ws.Range("D16").Value = "'65.94"
$ws.Range("D16").Style = $origStyle
$ws.Range("E16").Value = "  +1.50%  "
$ws.Range("D17").Value = "26.541.42"
$ws.Range("E18").Value = "  -0.02%  "
$origStyle = $ws.Range("D19").Style
$ws.Range("D19").Value = "'4.804"
$ws.Range("D19").Style = $origStyle
$ws.Range("E19").Value = "  +3.51%  "
$origStyle = $ws.Range("D20").Style
$ws.Range("D20").Value = "'195.24"
$ws.Range("D20").Style = $origStyle
$ws.Range("E20").Value = "  +5.34%  "
$origStyle = $ws.Range("D21").Style
$ws.Range("D21").Value = "'10.41"
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = "  +3.79%  "
$origStyle = $ws.Range("D22").Style
$ws.Range("D22").Value = "'6.387"
$ws.Range("D22").Style = $origStyle
$ws.Range("E22").Value = "  +4.94%  "
$origStyle = $ws.Range("D23").Style
$ws.Range("D23").Value = "'1.002"
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = "  -0.09%  "
$origStyle = $ws.Range("D24").Style
$ws.Range("D24").Value = "'143.22"
$ws.Range("D24").Style = $origStyle
$ws.Range("E24").Value = "  -1.47%  "
$origStyle = $ws.Range("D25").Style
$ws.Range("D25").Value = "'0.1276"
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = "  +7.12%  "
$origStyle = $ws.Range("D26").Style
$ws.Range("D26").Value = "'7.471"
$ws.Range("D26").Style = $origStyle
$ws.Range("E26").Value = "  +1.65%  "
$ws.Range("E27").Value = "  +4.57%  "
$origStyle = $ws.Range("D28").Style
$ws.Range("D28").Value = "'1.414"
$ws.Range("D28").Style = $origStyle
$ws.Range("E28").Value = "  +3.32%  "
$origStyle = $ws.Range("D29").Style
$ws.Range("D29").Value = "'0.06131"
$ws.Range("D29").Style = $origStyle
$ws.Range("E29").Value = "  +2.61%  "
$origStyle = $ws.Range("D30").Style
$ws.Range("D30").Value = "'1.278"
$ws.Range("D30").Style = $origStyle
$ws.Range("E30").Value = "  +2.68%  "
$origStyle = $ws.Range("D31").Style
$ws.Range("D31").Value = "'3.608"
$ws.Range("D31").Style = $origStyle
$ws.Range("E31").Value = "  +7.58%  "
$origStyle = $ws.Range("D32").Style
$ws.Range("D32").Value = "'3.459"
$ws.Range("D32").Style = $origStyle
$ws.Range("E32").Value = "  +3.63%  "
$origStyle = $ws.Range("D33").Style
$ws.Range("D33").Value = "'1.706"
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = "  +6.08%  "
$origStyle = $ws.Range("D34").Style
$ws.Range("D34").Value = "'1.015"
$ws.Range("D34").Style = $origStyle
$ws.Range("E34").Value = "  +4.65%  "
$origStyle = $ws.Range("D35").Style
$ws.Range("D35").Value = "'2.791"
$ws.Range("D35").Style = $origStyle
$ws.Range("E35").Value = "  +2.33%  "
$ws.Range("E36").Value = "  +1.62%  "
$origStyle = $ws.Range("D37").Style
$ws.Range("D37").Value = "'0.5706"
$ws.Range("D37").Style = $origStyle
$ws.Range("E37").Value = "  -2.18%  "
$origStyle = $ws.Range("D38").Style
$ws.Range("D38").Value = "'0.01642"
$ws.Range("D38").Style = $origStyle
$ws.Range("E38").Value = "  +3.13%  "
$origStyle = $ws.Range("D39").Style
$ws.Range("D39").Value = "'5.949"
$ws.Range("D39").Style = $origStyle
$ws.Range("E39").Value = "  +3.31%  "
$origStyle = $ws.Range("D40").Style
$ws.Range("D40").Value = "'0.8711"
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = "  +3.70%  "
$ws.Range("D41").Value = "1.061.96"
$ws.Range("E41").Value = "  +1.41%  "
$ws.Range("E42").Value = "  -0.07%  "
$origStyle = $ws.Range("D43").Style
$ws.Range("D43").Value = "'100.00"
$ws.Range("D43").Style = $origStyle
$ws.Range("E43").Value = "  +0.28%  "
$ws.Range("D44").Value = "1.831.75"
$ws.Range("D45").Value = "0.0₈112"
$ws.Range("E45").Value = "  +4.58%  "
$origStyle = $ws.Range("D46").Style
$ws.Range("D46").Value = "'57.21"
$ws.Range("D46").Style = $origStyle
$ws.Range("E46").Value = "  +5.65%  "
$origStyle = $ws.Range("D47").Style
$ws.Range("D47").Value = "'8.147"
$ws.Range("D47").Style = $origStyle
$ws.Range("E47").Value = "  +2.20%  "
$origStyle = $ws.Range("D48").Style
$ws.Range("D48").Value = "'0.9981"
$ws.Range("D48").Style = $origStyle
$ws.Range("E48").Value = "  -0.20%  "
$origStyle = $ws.Range("D49").Style
$ws.Range("D49").Value = "'0.05204"
$ws.Range("D49").Style = $origStyle
$ws.Range("E49").Value = "  +0.24%  "
$origStyle = $ws.Range("D50").Style
$ws.Range("D50").Value = "'6.073"
$ws.Range("D50").Style = $origStyle
$ws.Range("E50").Value = "  +4.74%  "
$origStyle = $ws.Range("D51").Style
$ws.Range("D51").Value = "'0.4241"
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = "  +0.23%  "
